$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be parsed as numbers
$textCells = @("D5", "D6", "D9", "D11", "D12", "D14", "D16", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D29", "D30", "D31", "D32", "D34", "D35", "D36", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D50", "D51")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

# Apply updated values
$ws.Range("D2").Value = '61.031.45'
$ws.Range("E2").Value = '  -0.88%  '
$ws.Range("D3").Value = '3.414.11'
$ws.Range("E3").Value = '  -1.03%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").Value = '572.56'
$ws.Range("E5").Value = '  -1.22%  '
$ws.Range("D6").Value = '143.55'
$ws.Range("E6").Value = '  -2.79%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = '3.415.80'
$ws.Range("E8").Value = '  -0.97%  '
$ws.Range("D9").Value = '0.476'
$ws.Range("E9").Value = '  +0.83%  '
$ws.Range("E10").Value = '  -0.82%  '
$ws.Range("D11").Value = '0.126'
$ws.Range("E11").Value = '  +1.78%  '
$ws.Range("D12").Value = '0.393'
$ws.Range("E12").Value = '  +1.38%  '
$ws.Range("D13").Value = '3.995.49'
$ws.Range("E13").Value = '  -1.07%  '
$ws.Range("D14").Value = '28.39'
$ws.Range("E14").Value = '  +2.05%  '
$ws.Range("E15").Value = '  +0.65%  '
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").Value = '0.0000172'
$ws.Range("E16").Value = '  -1.14%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.414.92'
$ws.Range("E17").Value = '  -1.00%  '
$ws.Range("D18").Value = '61.107.00'
$ws.Range("E18").Value = '  -0.93%  '
$ws.Range("D19").Value = '6.32'
$ws.Range("E19").Value = '  +0.35%  '
$ws.Range("D20").Value = '14.44'
$ws.Range("E20").Value = '  +2.56%  '
$ws.Range("D21").Value = '9.39'
$ws.Range("E21").Value = '  -0.72%  '
$ws.Range("D22").Value = '389.90'
$ws.Range("E22").Value = '  +1.39%  '
$ws.Range("D23").Value = '0.570'
$ws.Range("E23").Value = '  +0.81%  '
$ws.Range("D24").Value = '72.83'
$ws.Range("E24").Value = '  +0.82%  '
$ws.Range("D25").Value = '0.996'
$ws.Range("E25").Value = '  -0.38%  '
$ws.Range("D26").Value = '0.0000124'
$ws.Range("E26").Value = '  +0.51%  '
$ws.Range("D27").Value = '3.553.43'
$ws.Range("E27").Value = '  -1.00%  '
$ws.Range("E28").Value = '  +1.96%  '
$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D29").Value = '7.49'
$ws.Range("E29").Value = '  -3.58%  '
$ws.Range("B30").Value = 'Binance-PegBSC-USD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D30").Value = '0.998'
$ws.Range("E30").Value = '  -0.20%  '
$ws.Range("D31").Value = '8.21'
$ws.Range("E31").Value = '  +0.15%  '
$ws.Range("D32").Value = '1.46'
$ws.Range("E32").Value = '  -6.22%  '
$ws.Range("E33").Value = '  -0.06%  '
$ws.Range("D34").Value = '0.999'
$ws.Range("E34").Value = '  -0.15%  '
$ws.Range("D35").Value = '23.96'
$ws.Range("E35").Value = '  -0.20%  '
$ws.Range("D36").Value = '7.05'
$ws.Range("E36").Value = '  +0.38%  '
$ws.Range("D37").Value = '3.441.09'
$ws.Range("E37").Value = '  -0.85%  '
$ws.Range("D38").Value = '5.16'
$ws.Range("E38").Value = '  -1.42%  '
$ws.Range("B39").Value = 'ImmutableX'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D39").Value = '1.55'
$ws.Range("E39").Value = '  -0.49%  '
$ws.Range("B40").Value = 'Monero'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D40").Value = '167.83'
$ws.Range("E40").Value = '  +1.28%  '
$ws.Range("D41").Value = '0.0788'
$ws.Range("E41").Value = '  +0.35%  '
$ws.Range("D42").Value = '27.21'
$ws.Range("E42").Value = '  +3.99%  '
$ws.Range("D43").Value = '0.795'
$ws.Range("E43").Value = '  +0.40%  '
$ws.Range("B44").Value = 'Filecoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D44").Value = '4.51'
$ws.Range("E44").Value = '  +0.81%  '
$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").Value = '0.999'
$ws.Range("E45").Value = '  -0.17%  '
$ws.Range("D46").Value = '41.99'
$ws.Range("E46").Value = '  +0.00%  '
$ws.Range("D47").Value = '1.72'
$ws.Range("E47").Value = '  -0.77%  '
$ws.Range("D48").Value = '2.588.06'
$ws.Range("E48").Value = '  -1.08%  '
$ws.Range("E49").Value = '  -3.04%  '
$ws.Range("D50").Value = '6.98'
$ws.Range("E50").Value = '  +1.61%  '
$ws.Range("D51").Value = '23.08'
$ws.Range("E51").Value = '  -2.57%  '
